# Edit the "Invoice 1.xlsx" workbook:
#  1. Update the CLIENT CODE placeholder footer text with an actual name/email.
#  2. Update the "Client discount" unit price and format it as currency.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Invoice")

# 1) Footer placeholder text -> actual client contact info
$ws.Range("A31").Value = "Charlie, charlie@mail.com"

# 2) Client discount unit price (column E, row 18) -> 100, formatted as currency
$ws.Range("E18").Value = 100
$ws.Range("E18").NumberFormat = '"$"#,##0_);[Red]\("$"#,##0\)'
